$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.830.25'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '1.643.55'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.84'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.506'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.253'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.71'
$ws.Range('E10').Value = '  +3.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '1.873.01'
$ws.Range('D13').Value = '1.640.49'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.16'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').Value = '26.854.99'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.59'
$ws.Range('E19').Value = '  +2.97%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('E22').Value = '  +7.12%  '
$ws.Range('E23').Value = '  +7.20%  '
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.02'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.36'
$ws.Range('E27').Value = '  +3.94%  '
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('E29').Value = '  +1.68%  '
$ws.Range('E30').Value = '  +1.72%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').Value = '1.248.26'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0174'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.535'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.834'
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.806'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').Value = '1.785.15'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('E44').Value = '  -3.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.81'
$ws.Range('E45').Value = '  +1.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.50'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0971'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.54'
$ws.Range('E51').Value = '  +0.31%  '
